$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Make "Repayment schedule" the active sheet/tab (was "Transactions" before).
$ws.Activate()

# Insert a new (blank) column before column N, shifting
# N(Late) -> O, O(heading/Outstanding) -> P, P(Outstanding) -> Q.
$ws.Columns("N").Insert()

# The freshly inserted column picks up the width of its left neighbour
# (column M, "In Advance") rather than an auto best-fit width.
$ws.Columns("N").ColumnWidth = 9.83

# Update the selection shown on the Repayment schedule sheet.
$ws.Range("K13").Select()
